$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Calificaciones sheet: fill in the previously-missing 2P grades for
#    "Conciencia historica" (col K) and "Reacciones quimicas" (col P), and
#    update the recalculated "Final" grades for those same two subjects
#    (cols AA and AF) for every student row (4-27).
# ---------------------------------------------------------------------------
$wsCal = $wb.Worksheets.Item("Calificaciones")

$grades = @{
    4  = @{ K = 10; P = 9;  AA = 9  }
    5  = @{ K = 10; P = 7;  AA = 9;  AF = 8 }
    6  = @{ K = 10; P = 9;  AA = 10 }
    7  = @{ K = 7;  P = 6;  AA = 6  }
    8  = @{ K = 10; P = 10 }
    9  = @{ K = 7;  P = 8;  AA = 7;  AF = 8 }
    10 = @{ K = 9;  P = 8;  AA = 8  }
    11 = @{ K = 7;  P = 9;  AF = 8 }
    12 = @{ K = 10; P = 9;  AA = 9  }
    13 = @{ K = 10; P = 9;  AA = 9  }
    14 = @{ K = 10; P = 10 }
    15 = @{ K = 10; P = 7;  AA = 8  }
    16 = @{ K = 5;  P = 6;  AA = 5  }
    17 = @{ K = 6;  P = 6  }
    18 = @{ K = 5;  P = 6  }
    19 = @{ K = 8;  P = 8  }
    20 = @{ K = 5;  P = 8;  AF = 7 }
    21 = @{ K = 5;  P = 6;  AA = 5;  AF = 7 }
    22 = @{ K = 10; P = 8  }
    23 = @{ K = 6;  P = 7  }
    24 = @{ K = 7;  P = 6;  AA = 7  }
    25 = @{ K = 7;  P = 6;  AA = 6  }
    26 = @{ K = 6;  P = 8  }
    27 = @{ K = 6;  P = 6  }
}

foreach ($row in $grades.Keys) {
    $cols = $grades[$row]
    foreach ($col in $cols.Keys) {
        $wsCal.Range("$col$row").Value = $cols[$col]
    }
}

# ---------------------------------------------------------------------------
# 2) Totales sheet: the "Conciencia historica" average (H2) is recalculated
#    now that every student has a 2P grade for that subject.
# ---------------------------------------------------------------------------
$wsTot = $wb.Worksheets.Item("Totales")
$wsTot.Range("H2").Value = 7.3

# ---------------------------------------------------------------------------
# 3) Rescatables sheet: refreshed list of students who still need to pass a
#    pending subject, with updated control numbers / names / subjects /
#    teachers for the second partial.
# ---------------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Rescatables")

$rescatables = @(
    @{ Row = 2;  NC = 23330051920018; Paterno = "RAMOS";       Materno = "UTRERA";     Nombres = "CARLOS DAVID"; Materia = "Conciencia histórica I. Perspectivas del México antiguo en los contextos globales"; Docente = "Vargas Olvera Francisco Eduardo" }
    @{ Row = 3;  NC = 23330051920018; Paterno = "RAMOS";       Materno = "UTRERA";     Nombres = "CARLOS DAVID"; Materia = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"; Docente = "Jiménez Nieto Enrique" }
    @{ Row = 4;  NC = 23330051920018; Paterno = "RAMOS";       Materno = "UTRERA";     Nombres = "CARLOS DAVID"; Materia = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"; Docente = "Jiménez Nieto Enrique" }
    @{ Row = 5;  NC = 23330051920005; Paterno = "CASTRO";      Materno = "ARIAS";      Nombres = "OMAR DAVID";   Materia = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"; Docente = "Jiménez Nieto Enrique" }
    @{ Row = 6;  NC = 23330051920005; Paterno = "CASTRO";      Materno = "ARIAS";      Nombres = "OMAR DAVID";   Materia = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"; Docente = "Jiménez Nieto Enrique" }
    @{ Row = 7;  NC = 23330051920025; Paterno = "XOTLANIHUA";  Materno = "ZEPAHUA";    Nombres = "JUAN ALBERTO"; Materia = "Conciencia histórica I. Perspectivas del México antiguo en los contextos globales"; Docente = "Vargas Olvera Francisco Eduardo" }
    @{ Row = 8;  NC = 23330051920025; Paterno = "XOTLANIHUA";  Materno = "ZEPAHUA";    Nombres = "JUAN ALBERTO"; Materia = "Reacciones químicas, conservación de la materia en la formación de nuevas substancias."; Docente = "Nativitas Sandoval Liliana Soledad" }
    @{ Row = 9;  NC = 23330051920014; Paterno = "JUSTO";       Materno = "NEGRETE";    Nombres = "JAQUELINE";    Materia = "Conciencia histórica I. Perspectivas del México antiguo en los contextos globales"; Docente = "Vargas Olvera Francisco Eduardo" }
    @{ Row = 10; NC = 23330051920321; Paterno = "ROMERO";      Materno = "RAMOS";      Nombres = "LUIS DAVID";   Materia = "Conciencia histórica I. Perspectivas del México antiguo en los contextos globales"; Docente = "Vargas Olvera Francisco Eduardo" }
    @{ Row = 11; NC = 23330051920023; Paterno = "VASQUEZ";     Materno = "ESPINDOLA";  Nombres = "JOSUE YAHIR";  Materia = "Reacciones químicas, conservación de la materia en la formación de nuevas substancias."; Docente = "Nativitas Sandoval Liliana Soledad" }
)

foreach ($item in $rescatables) {
    $r = $item.Row
    $wsRes.Range("A$r").Value = $item.NC
    $wsRes.Range("B$r").Value = $item.Paterno
    $wsRes.Range("C$r").Value = $item.Materno
    $wsRes.Range("D$r").Value = $item.Nombres
    $wsRes.Range("E$r").Value = $item.Materia
    $wsRes.Range("F$r").Value = $item.Docente
    $wsRes.Range("G$r").Value = 5
}
